$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 3346.3845
$ws.Range("I74").Value = 3643.2856
$ws.Range("J74").Value = 3000
$ws.Range("K74").Value = 3643.2856
$ws.Range("L74").Value = 3000
$ws.Range("M74").Value = -2707.2856
$ws.Range("N74").Value = -4872

$ws.Range("H76").Value = 3014.2856
$ws.Range("I76").Value = 2936.3635
$ws.Range("J76").Value = 3300
$ws.Range("K76").Value = 2936.3635
$ws.Range("L76").Value = 3300
$ws.Range("M76").Value = -2621.3635
$ws.Range("N76").Value = -3930

$ws.Range("H77").Value = 3346.3845
$ws.Range("I77").Value = 3643.2856
$ws.Range("J77").Value = 3000
$ws.Range("K77").Value = 18216.428
$ws.Range("L77").Value = 15000
$ws.Range("M77").Value = -13536.428
$ws.Range("N77").Value = -24360

$ws.Range("H79").Value = 3014.2856
$ws.Range("I79").Value = 2936.3635
$ws.Range("J79").Value = 3300
$ws.Range("K79").Value = 2936.3635
$ws.Range("L79").Value = 3300
$ws.Range("M79").Value = -1844.3635
$ws.Range("N79").Value = -5484

$ws.Range("H129").Value = 1061.7812
$ws.Range("J129").Value = 1061.7812
$ws.Range("L129").Value = 3185.3436
$ws.Range("N129").Value = -13185.3436

$ws.Range("H135").Value = 1202225.1
$ws.Range("I135").Value = 578.4666999999999
$ws.Range("J135").Value = 3004695.2
$ws.Range("K135").Value = 5206.2003
$ws.Range("L135").Value = 27042256.8
$ws.Range("M135").Value = -2671.2003
$ws.Range("N135").Value = -27047326.8

$ws.Range("H141").Value = 1583.1428
$ws.Range("I141").Value = 799.7143
$ws.Range("K141").Value = 2399.1429
$ws.Range("M141").Value = 2780.8571

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H55").Value = 17933.666
$ws.Range("J55").Value = 17933.666
$ws.Range("L55").Value = 17933.666
$ws.Range("N55").Value = -18563.666

$ws.Range("H63").Value = 3796.6667
$ws.Range("I63").Value = 2595.7144
$ws.Range("J63").Value = 8000
$ws.Range("K63").Value = 2595.7144
$ws.Range("L63").Value = 8000
$ws.Range("M63").Value = -1909.7144
$ws.Range("N63").Value = -9372

$ws.Range("H66").Value = 3796.6667
$ws.Range("I66").Value = 2595.7144
$ws.Range("J66").Value = 8000
$ws.Range("K66").Value = 12978.572
$ws.Range("L66").Value = 40000
$ws.Range("M66").Value = -9546.572
$ws.Range("N66").Value = -46864

$ws.Range("H102").Value = 2013.75
$ws.Range("I102").Value = 2013.75
$ws.Range("K102").Value = 2013.75
$ws.Range("M102").Value = -391.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2172
$ws.Range("I105").Value = 2362.8
$ws.Range("K105").Value = 2362.8
$ws.Range("M105").Value = -615.8000000000002

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H50").Value = 10271.857
$ws.Range("J50").Value = 10483.833
$ws.Range("L50").Value = 10483.833
$ws.Range("N50").Value = -11733.833

$ws.Range("H51").Value = 10428.714
$ws.Range("J51").Value = 10650.167
$ws.Range("L51").Value = 10650.167
$ws.Range("N51").Value = -12122.167

$ws.Range("H59").Value = 16444.111
$ws.Range("J59").Value = 16444.111
$ws.Range("L59").Value = 16444.111
$ws.Range("N59").Value = -18734.111

$ws.Range("H60").Value = 9330.6
$ws.Range("J60").Value = 10163.25
$ws.Range("L60").Value = 10163.25
$ws.Range("N60").Value = -11185.25

$ws.Range("H61").Value = 10428.714
$ws.Range("J61").Value = 10650.167
$ws.Range("L61").Value = 10650.167
$ws.Range("N61").Value = -11346.167

$ws.Range("H68").Value = 20239
$ws.Range("J68").Value = 20239
$ws.Range("L68").Value = 20239
$ws.Range("N68").Value = -21737

$ws.Range("H71").Value = 20239
$ws.Range("J71").Value = 20239
$ws.Range("L71").Value = 60717
$ws.Range("N71").Value = -68205

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 4235.3667
$ws.Range("I131").Value = 8743.5
$ws.Range("J131").Value = 3108.3333
$ws.Range("K131").Value = 26230.5
$ws.Range("L131").Value = 9324.999899999999
$ws.Range("M131").Value = -21190.5
$ws.Range("N131").Value = -19404.9999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4712.5
$ws.Range("I70").Value = 4875
$ws.Range("J70").Value = 4550
$ws.Range("K70").Value = 4875
$ws.Range("L70").Value = 4550
$ws.Range("M70").Value = -4605
$ws.Range("N70").Value = -5090

$ws.Range("H73").Value = 4712.5
$ws.Range("I73").Value = 4875
$ws.Range("J73").Value = 4550
$ws.Range("K73").Value = 4875
$ws.Range("L73").Value = 4550
$ws.Range("M73").Value = -3939
$ws.Range("N73").Value = -6422

$ws.Range("H80").Value = 6211.8965
$ws.Range("I80").Value = 2436.3333
$ws.Range("J80").Value = 10257.143
$ws.Range("K80").Value = 2436.3333
$ws.Range("L80").Value = 10257.143
$ws.Range("M80").Value = -1438.3333
$ws.Range("N80").Value = -12253.143

$ws.Range("H83").Value = 6211.8965
$ws.Range("I83").Value = 2436.3333
$ws.Range("J83").Value = 10257.143
$ws.Range("K83").Value = 12181.6665
$ws.Range("L83").Value = 51285.715
$ws.Range("M83").Value = -7189.666499999999
$ws.Range("N83").Value = -61269.715

$ws.Range("H97").Value = 58824940

$ws.Range("H126").Value = 2422.8
$ws.Range("I126").Value = 2333.3333
$ws.Range("J126").Value = 2557
$ws.Range("K126").Value = 6999.999899999999
$ws.Range("L126").Value = 7671
$ws.Range("M126").Value = -4529.999899999999
$ws.Range("N126").Value = -12611

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 2447.3076
$ws.Range("I93").Value = 2985.1667
$ws.Range("J93").Value = 1986.2858
$ws.Range("K93").Value = 2985.1667
$ws.Range("L93").Value = 1986.2858
$ws.Range("M93").Value = -1737.1667
$ws.Range("N93").Value = -4482.2858

$ws.Range("H132").Value = 4659.943
$ws.Range("I132").Value = 4930.923
$ws.Range("J132").Value = 3877.111
$ws.Range("K132").Value = 14792.769
$ws.Range("L132").Value = 11631.333
$ws.Range("M132").Value = -12262.769
$ws.Range("N132").Value = -16691.333

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 391.38095
$ws.Range("I107").Value = 204
$ws.Range("J107").Value = 695.875
$ws.Range("K107").Value = 612
$ws.Range("L107").Value = 2087.625
$ws.Range("M107").Value = 1308
$ws.Range("N107").Value = -5927.625
